# Updated capital structure database
# Applies diff: row 2 values updated, new row 3 (Banco Pichincha) inserted,
# former row 3 (Interbank) data shifted to row 4 with refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range('A2').Value = 'Peru'
$ws.Range('B2').Value = "'2"
$ws.Range('C2').Value = 'Banks (Regional)'
$ws.Range('D2').Value = -0.01514
$ws.Range('E2').Value = -0.1422
$ws.Range('G2').Value = 0
$ws.Range('H2').Value = 0
$ws.Range('I2').Value = 0
$ws.Range('J2').Value = 0
$ws.Range('K2').Value = 149.03
$ws.Range('L2').Value = 0.2104646236407287
$ws.Range('M2').Value = 84.1
$ws.Range('N2').Value = 0.01681058607179979
$ws.Range('O2').Value = 0.564315909548413
$ws.Range('P2').Value = 84.1
$ws.Range('Q2').Value = 0.01681058607179979
$ws.Range('R2').Value = 0.564315909548413
$ws.Range('S2').Value = 0
$ws.Range('T2').Value = 0
$ws.Range('U2').Value = 2533.5
$ws.Range('V2').Value = 0.5064164068121851
$ws.Range('W2').Value = 0.05429876105795206
$ws.Range('X2').Value = 0.1164721598537399
$ws.Range('Y2').Value = -0.06217339879578788
$ws.Range('Z2').Value = 0.1356833123513778
$ws.Range('AA2').Value = 0
$ws.Range('AB2').Value = 0.04396422323290067
$ws.Range('AC2').Value = -0.04396422323290067
$ws.Range('AD2').Value = 5453.6
$ws.Range('AE2').Value = 0
$ws.Range('AF2').Value = 5453.6
$ws.Range('AG2').Value = 2920.1
$ws.Range('AH2').Value = 0.5215561761218009
$ws.Range('AI2').Value = 0.7371722087050555
$ws.Range('AJ2').Value = 0.3685645407615898
$ws.Range('AK2').Value = 0.60028779936273
$ws.Range('AL2').Value = 0
$ws.Range('AM2').Value = 0

# ---- Row 3 ----
$ws.Range('A3').Value = 'Peru'
$ws.Range('B3').Value = 'Banco Pichincha S.A. (BVL:BPICHC1)'
$ws.Range('C3').Value = 'Banks (Regional)'
$ws.Range('D3').Value = 0.00212
$ws.Range('E3').Value = -0.186
$ws.Range('G3').Value = 0
$ws.Range('H3').Value = 0
$ws.Range('I3').Value = 0
$ws.Range('J3').Value = 0
$ws.Range('K3').Value = 6.83
$ws.Range('L3').Value = 0.05965065502183406
$ws.Range('M3').Value = -0.0
$ws.Range('N3').Value = -0.0
$ws.Range('O3').Value = -0.0
$ws.Range('P3').Value = -0.0
$ws.Range('Q3').Value = -0.0
$ws.Range('R3').Value = -0.0
$ws.Range('S3').Value = 0
$ws.Range('T3').ClearContents()
$ws.Range('U3').Value = 454.9
$ws.Range('V3').Value = 3.713469387755102
$ws.Range('W3').Value = 0.02559010865492694
$ws.Range('X3').Value = 0.1699205629205485
$ws.Range('Y3').Value = -0.1443304542656215
$ws.Range('Z3').Value = 0.231126362535325
$ws.Range('AA3').Value = 0
$ws.Range('AB3').Value = 0.04312754993837023
$ws.Range('AC3').Value = -0.04312754993837023
$ws.Range('AD3').Value = 702.3
$ws.Range('AE3').Value = 0
$ws.Range('AF3').Value = 702.3
$ws.Range('AG3').Value = 247.4
$ws.Range('AH3').Value = 0.8514791464597478
$ws.Range('AI3').Value = 0.7295107510127765
$ws.Range('AJ3').Value = 0.6688294133549608
$ws.Range('AK3').Value = 0.487199684915321
$ws.Range('AL3').Value = 0
$ws.Range('AM3').Value = 0

# ---- Row 4 ----
$ws.Range('A4').Value = 'Peru'
$ws.Range('B4').Value = 'Banco Internacional del Perú S.A.A. - Interbank (BVL:INTERBC1)'
$ws.Range('C4').Value = 'Banks (Regional)'
$ws.Range('D4').Value = -0.03240000000000001
$ws.Range('E4').Value = -0.0984
$ws.Range('G4').Value = 0
$ws.Range('H4').Value = 0
$ws.Range('I4').Value = 0
$ws.Range('J4').Value = 0
$ws.Range('K4').Value = 142.2
$ws.Range('L4').Value = 0.23955525606469
$ws.Range('M4').Value = 84.1
$ws.Range('N4').Value = 0.01723254717947667
$ws.Range('O4').Value = 0.5914205344585092
$ws.Range('P4').Value = 84.1
$ws.Range('Q4').Value = 0.01723254717947667
$ws.Range('R4').Value = 0.5914205344585092
$ws.Range('S4').Value = 0
$ws.Range('T4').Value = 0
$ws.Range('U4').Value = 2078.6
$ws.Range('V4').Value = 0.4259164395631416
$ws.Range('W4').Value = 0.08300741346097718
$ws.Range('X4').Value = 0.0630237567869314
$ws.Range('Y4').Value = 0.01998365667404578
$ws.Range('Z4').Value = 0.1256729834842496
$ws.Range('AA4').Value = 0
$ws.Range('AB4').Value = 0.04480089652743111
$ws.Range('AC4').Value = -0.04480089652743111
$ws.Range('AD4').Value = 4751.3
$ws.Range('AE4').Value = 0
$ws.Range('AF4').Value = 4751.3
$ws.Range('AG4').Value = 2672.7
$ws.Range('AH4').Value = 0.4933032933261348
$ws.Range('AI4').Value = 0.7383183379174242
$ws.Range('AJ4').Value = 0.3538593936184298
$ws.Range('AK4').Value = 0.6134689099547822
$ws.Range('AL4').Value = 0
$ws.Range('AM4').Value = 0
